$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws 'D2' '57.870.82'
Set-TextValue $ws 'D3' '3.116.50'
Set-TextValue $ws 'E3' '  -1.00%  '
Set-TextValue $ws 'E4' '  +0.05%  '
Set-TextValue $ws 'D5' '529.90'
Set-TextValue $ws 'E5' '  -0.38%  '
Set-TextValue $ws 'D6' '138.02'
Set-TextValue $ws 'D7' '0.999'
Set-TextValue $ws 'E7' '  -0.17%  '
Set-TextValue $ws 'E8' '  +3.49%  '
Set-TextValue $ws 'D9' '7.28'
Set-TextValue $ws 'E9' '  +0.79%  '
Set-TextValue $ws 'E10' '  -2.22%  '
Set-TextValue $ws 'E11' '  +1.87%  '
Set-TextValue $ws 'D12' '3.659.38'
Set-TextValue $ws 'E12' '  -0.85%  '
Set-TextValue $ws 'E13' '  +1.35%  '
Set-TextValue $ws 'D14' '25.42'
Set-TextValue $ws 'E14' '  -0.94%  '
Set-TextValue $ws 'E15' '  -2.16%  '
Set-TextValue $ws 'D16' '57.909.88'
Set-TextValue $ws 'E16' '  -0.01%  '
Set-TextValue $ws 'D17' '3.121.40'
Set-TextValue $ws 'E17' '  -0.81%  '
Set-TextValue $ws 'E18' '  -2.24%  '
Set-TextValue $ws 'D19' '12.56'
Set-TextValue $ws 'E19' '  -2.04%  '
Set-TextValue $ws 'D20' '7.98'
Set-TextValue $ws 'E20' '  -0.20%  '
Set-TextValue $ws 'D21' '350.45'
Set-TextValue $ws 'E21' '  -1.32%  '
Set-TextValue $ws 'E22' '  -0.15%  '
Set-TextValue $ws 'D23' '68.96'
Set-TextValue $ws 'E23' '  +0.62%  '
Set-TextValue $ws 'D24' '0.503'
Set-TextValue $ws 'E24' '  -1.35%  '
Set-TextValue $ws 'E25' '  -1.59%  '
Set-TextValue $ws 'D26' '0.998'
Set-TextValue $ws 'E26' '  -0.38%  '
Set-TextValue $ws 'D27' '0.0₃0869'
Set-TextValue $ws 'E27' '  -7.45%  '
Set-TextValue $ws 'D28' '7.20'
Set-TextValue $ws 'E28' '  -3.49%  '
Set-TextValue $ws 'D29' '1.86'
Set-TextValue $ws 'E29' '  -2.00%  '
Set-TextValue $ws 'E30' '  -5.34%  '
Set-TextValue $ws 'D31' '21.18'
Set-TextValue $ws 'E31' '  -0.48%  '
Set-TextValue $ws 'D32' '4.91'
Set-TextValue $ws 'E32' '  +0.62%  '
Set-TextValue $ws 'E33' '  -5.44%  '
Set-TextValue $ws 'D34' '158.66'
Set-TextValue $ws 'E34' '  +0.80%  '
Set-TextValue $ws 'D35' '6.03'
Set-TextValue $ws 'E35' '  -2.76%  '
Set-TextValue $ws 'D36' '25.83'
Set-TextValue $ws 'E36' '  -1.23%  '
Set-TextValue $ws 'E37' '  -2.45%  '
Set-TextValue $ws 'D38' '1.66'
Set-TextValue $ws 'E38' '  +1.93%  '
Set-TextValue $ws 'D39' '0.0668'
Set-TextValue $ws 'E39' '  -0.64%  '
Set-TextValue $ws 'D40' '4.00'
Set-TextValue $ws 'E40' '  -1.80%  '
Set-TextValue $ws 'D41' '0.695'
Set-TextValue $ws 'E41' '  -1.07%  '
Set-TextValue $ws 'D42' '37.14'
Set-TextValue $ws 'E42' '  +1.19%  '
Set-TextValue $ws 'D43' '2.391.15'
Set-TextValue $ws 'E43' '  +2.61%  '
Set-TextValue $ws 'D44' '3.161.46'
Set-TextValue $ws 'E44' '  -0.90%  '
Set-TextValue $ws 'D45' '1.00'
Set-TextValue $ws 'E45' '  +0.11%  '
Set-TextValue $ws 'D46' '0.0265'
Set-TextValue $ws 'E46' '  -2.94%  '
Set-TextValue $ws 'D47' '0.960'
Set-TextValue $ws 'E47' '  -4.29%  '
Set-TextValue $ws 'E48' '  -0.68%  '
Set-TextValue $ws 'D49' '19.70'
Set-TextValue $ws 'E49' '  -3.19%  '
Set-TextValue $ws 'D50' '0.737'
Set-TextValue $ws 'E50' '  -2.28%  '
Set-TextValue $ws 'D51' '0.0911'
Set-TextValue $ws 'E51' '  +1.71%  '
